$d = $word.ActiveDocument

# 1. Insert the new checklist item "Complete one of the learning path items"
#    right before "Work somewhere other than your desk" (currently paragraph 9).
$anchor = $d.Paragraphs.Item(9)
$anchor.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item(9)
$newPara.Range.Text = "Complete one of the learning path items"

# 2. Update the line spacing of every checklist ("ListParagraph") item from
#    480/auto to 640/exact (i.e. exactly 32pt).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "List Paragraph") {
        $p.Format.LineSpacingRule = 4
        $p.Format.LineSpacing = 32
    }
}
